$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.65"
$ws.Range("E2").Value = "'-1.04%"
$ws.Range("D3").Value = "'37.43"
$ws.Range("E3").Value = "'-0.68%"
$ws.Range("D5").Value = "'0.07806"
$ws.Range("E5").Value = "'0.43%"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "'8.245"
$ws.Range("E6").Value = "'0.22%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.881"
$ws.Range("E7").Value = "'-0.12%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.997"
$ws.Range("E8").Value = "'2.96%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9271"
$ws.Range("E9").Value = "'0.68%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1085"
$ws.Range("E10").Value = "'-10.09%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1916"
$ws.Range("E11").Value = "'-0.55%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08933"
$ws.Range("E12").Value = "'-4.96%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03321"
$ws.Range("E13").Value = "'-2.36%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09582"
$ws.Range("E14").Value = "'-1.04%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001377"
$ws.Range("E15").Value = "'0.76%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005702"
$ws.Range("E16").Value = "'-1.61%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.563"
$ws.Range("E17").Value = "'0.33%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.431"
$ws.Range("E18").Value = "'1.82%"
$ws.Range("E19").Value = "'2.04%"
$ws.Range("E20").Value = "'18.96%"
$ws.Range("D21").Value = "'0.1275"
$ws.Range("E21").Value = "'-1.68%"
$ws.Range("D22").Value = "'0.2589"
$ws.Range("E22").Value = "'-0.05%"
$ws.Range("D23").Value = "'0.04388"
$ws.Range("E23").Value = "'1.09%"
$ws.Range("D24").Value = "'0.001201"
$ws.Range("E24").Value = "'-1.00%"
$ws.Range("D25").Value = "'0.004251"
$ws.Range("E25").Value = "'-0.13%"
$ws.Range("D26").Value = "'0.0001304"
$ws.Range("E26").Value = "'0.33%"
$ws.Range("D39").Value = "'0.02176"
$ws.Range("E39").Value = "'3.31%"
$ws.Range("D40").Value = "'0.05016"
$ws.Range("E40").Value = "'0.63%"
$ws.Range("D41").Value = "'0.007450"
$ws.Range("E41").Value = "'-2.72%"
$ws.Range("D42").Value = "'0.1346"
$ws.Range("E42").Value = "'0.23%"
$ws.Range("D43").Value = "'0.008674"
$ws.Range("E43").Value = "'-12.31%"
$ws.Range("D44").Value = "'0.002042"
$ws.Range("E44").Value = "'-0.96%"
$ws.Range("D45").Value = "'0.007979"
$ws.Range("E45").Value = "'-9.52%"
$ws.Range("D46").Value = "'0.00006560"
$ws.Range("E46").Value = "'-1.43%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.12%"
$ws.Range("D48").Value = "'0.002867"
$ws.Range("E48").Value = "'-2.35%"
$ws.Range("D49").Value = "'0.001002"
$ws.Range("E49").Value = "'-16.54%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.12%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.12%"
